# "Generate Report for Archive"
# The localization status report is regenerated: the shared "Status" value
# for the acfbde12-... entry moves on from "Ready for handoff" to
# "In Translation" everywhere it is shown (the Overview roll-up sheet as
# well as each per-locale detail sheet), and the now-shorter text lets the
# Status column shrink to fit.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: status is mirrored into one column per locale (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null

# --- Per-locale detail sheets: status lives in column C ("Status") ---
$locales = @("zh-cn", "de-de")
foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale)
    $ws.Range("C2").Value = $newStatus
    $ws.Columns.Item(3).AutoFit() | Out-Null
}
